$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 5620
$ws.Range("B2").Value = 45741
$ws.Range("A3").Value = 5570
$ws.Range("B3").Value = 45741.01041666666
$ws.Range("A4").Value = 5530
$ws.Range("B4").Value = 45741.02083333334
$ws.Range("A5").Value = 5490
$ws.Range("B5").Value = 45741.03125
$ws.Range("A6").Value = 5460
$ws.Range("B6").Value = 45741.04166666666
$ws.Range("A7").Value = 5440
$ws.Range("B7").Value = 45741.05208333334
$ws.Range("A8").Value = 5420
$ws.Range("B8").Value = 45741.0625
$ws.Range("A9").Value = 5410
$ws.Range("B9").Value = 45741.07291666666
$ws.Range("A10").Value = 5400
$ws.Range("B10").Value = 45741.08333333334
$ws.Range("A11").Value = 5400
$ws.Range("B11").Value = 45741.09375
$ws.Range("A12").Value = 5410
$ws.Range("B12").Value = 45741.10416666666
$ws.Range("A13").Value = 5420
$ws.Range("B13").Value = 45741.11458333334
$ws.Range("A14").Value = 5440
$ws.Range("B14").Value = 45741.125
$ws.Range("A15").Value = 5470
$ws.Range("B15").Value = 45741.13541666666
$ws.Range("A16").Value = 5500
$ws.Range("B16").Value = 45741.14583333334
$ws.Range("A17").Value = 5540
$ws.Range("B17").Value = 45741.15625
$ws.Range("A18").Value = 5590
$ws.Range("B18").Value = 45741.16666666666
$ws.Range("A19").Value = 5650
$ws.Range("B19").Value = 45741.17708333334
$ws.Range("A20").Value = 5730
$ws.Range("B20").Value = 45741.1875
$ws.Range("A21").Value = 5820
$ws.Range("B21").Value = 45741.19791666666
$ws.Range("A22").Value = 5930
$ws.Range("B22").Value = 45741.20833333334
$ws.Range("A23").Value = 6050
$ws.Range("B23").Value = 45741.21875
$ws.Range("A24").Value = 6190
$ws.Range("B24").Value = 45741.22916666666
$ws.Range("A25").Value = 6330
$ws.Range("B25").Value = 45741.23958333334
$ws.Range("A26").Value = 6470
$ws.Range("B26").Value = 45741.25
$ws.Range("A27").Value = 6600
$ws.Range("B27").Value = 45741.26041666666
$ws.Range("A28").Value = 6710
$ws.Range("B28").Value = 45741.27083333334
$ws.Range("A29").Value = 6800
$ws.Range("B29").Value = 45741.28125
$ws.Range("A30").Value = 6850
$ws.Range("B30").Value = 45741.29166666666
$ws.Range("A31").Value = 6870
$ws.Range("B31").Value = 45741.30208333334
$ws.Range("A32").Value = 6860
$ws.Range("B32").Value = 45741.3125
$ws.Range("A33").Value = 6810
$ws.Range("B33").Value = 45741.32291666666
$ws.Range("A34").Value = 6730
$ws.Range("B34").Value = 45741.33333333334
$ws.Range("A35").Value = 6620
$ws.Range("B35").Value = 45741.34375
$ws.Range("A36").Value = 6510
$ws.Range("B36").Value = 45741.35416666666
$ws.Range("A37").Value = 6390
$ws.Range("B37").Value = 45741.36458333334
$ws.Range("A38").Value = 6270
$ws.Range("B38").Value = 45741.375
$ws.Range("A39").Value = 6160
$ws.Range("B39").Value = 45741.38541666666
$ws.Range("A40").Value = 6060
$ws.Range("B40").Value = 45741.39583333334
$ws.Range("A41").Value = 5970
$ws.Range("B41").Value = 45741.40625
$ws.Range("A42").Value = 5900
$ws.Range("B42").Value = 45741.41666666666
$ws.Range("A43").Value = 5840
$ws.Range("B43").Value = 45741.42708333334
$ws.Range("A44").Value = 5790
$ws.Range("B44").Value = 45741.4375
$ws.Range("A45").Value = 5750
$ws.Range("B45").Value = 45741.44791666666
$ws.Range("A46").Value = 5710
$ws.Range("B46").Value = 45741.45833333334
$ws.Range("A47").Value = 5670
$ws.Range("B47").Value = 45741.46875
$ws.Range("A48").Value = 5640
$ws.Range("B48").Value = 45741.47916666666
$ws.Range("A49").Value = 5620
$ws.Range("B49").Value = 45741.48958333334
$ws.Range("A50").Value = 5600
$ws.Range("B50").Value = 45741.5
$ws.Range("A51").Value = 5600
$ws.Range("B51").Value = 45741.51041666666
$ws.Range("A52").Value = 5600
$ws.Range("B52").Value = 45741.52083333334
$ws.Range("A53").Value = 5620
$ws.Range("B53").Value = 45741.53125
$ws.Range("A54").Value = 5640
$ws.Range("B54").Value = 45741.54166666666
$ws.Range("A55").Value = 5670
$ws.Range("B55").Value = 45741.55208333334
$ws.Range("A56").Value = 5710
$ws.Range("B56").Value = 45741.5625
$ws.Range("A57").Value = 5750
$ws.Range("B57").Value = 45741.57291666666
$ws.Range("A58").Value = 5800
$ws.Range("B58").Value = 45741.58333333334
$ws.Range("A59").Value = 5850
$ws.Range("B59").Value = 45741.59375
$ws.Range("A60").Value = 5900
$ws.Range("B60").Value = 45741.60416666666
$ws.Range("A61").Value = 5970
$ws.Range("B61").Value = 45741.61458333334
$ws.Range("A62").Value = 6060
$ws.Range("B62").Value = 45741.625
$ws.Range("A63").Value = 6150
$ws.Range("B63").Value = 45741.63541666666
$ws.Range("A64").Value = 6250
$ws.Range("B64").Value = 45741.64583333334
$ws.Range("A65").Value = 6350
$ws.Range("B65").Value = 45741.65625
$ws.Range("A66").Value = 6470
$ws.Range("B66").Value = 45741.66666666666
$ws.Range("A67").Value = 6570
$ws.Range("B67").Value = 45741.67708333334
$ws.Range("A68").Value = 6680
$ws.Range("B68").Value = 45741.6875
$ws.Range("A69").Value = 6810
$ws.Range("B69").Value = 45741.69791666666
$ws.Range("A70").Value = 6910
$ws.Range("B70").Value = 45741.70833333334
$ws.Range("A71").Value = 7030
$ws.Range("B71").Value = 45741.71875
$ws.Range("A72").Value = 7150
$ws.Range("B72").Value = 45741.72916666666
$ws.Range("A73").Value = 7280
$ws.Range("B73").Value = 45741.73958333334
$ws.Range("A74").Value = 7410
$ws.Range("B74").Value = 45741.75
$ws.Range("A75").Value = 7510
$ws.Range("B75").Value = 45741.76041666666
$ws.Range("A76").Value = 7570
$ws.Range("B76").Value = 45741.77083333334
$ws.Range("A77").Value = 7580
$ws.Range("B77").Value = 45741.78125
$ws.Range("A78").Value = 7570
$ws.Range("B78").Value = 45741.79166666666
$ws.Range("A79").Value = 7540
$ws.Range("B79").Value = 45741.80208333334
$ws.Range("A80").Value = 7500
$ws.Range("B80").Value = 45741.8125
$ws.Range("A81").Value = 7430
$ws.Range("B81").Value = 45741.82291666666
$ws.Range("A82").Value = 7320
$ws.Range("B82").Value = 45741.83333333334
$ws.Range("A83").Value = 7200
$ws.Range("B83").Value = 45741.84375
$ws.Range("A84").Value = 7080
$ws.Range("B84").Value = 45741.85416666666
$ws.Range("A85").Value = 6950
$ws.Range("B85").Value = 45741.86458333334
$ws.Range("A86").Value = 6810
$ws.Range("B86").Value = 45741.875
$ws.Range("A87").Value = 6660
$ws.Range("B87").Value = 45741.88541666666
$ws.Range("A88").Value = 6530
$ws.Range("B88").Value = 45741.89583333334
$ws.Range("A89").Value = 6380
$ws.Range("B89").Value = 45741.90625
$ws.Range("A90").Value = 6260
$ws.Range("B90").Value = 45741.91666666666
$ws.Range("A91").Value = 6150
$ws.Range("B91").Value = 45741.92708333334
$ws.Range("A92").Value = 6020
$ws.Range("B92").Value = 45741.9375
$ws.Range("A93").Value = 5900
$ws.Range("B93").Value = 45741.94791666666
$ws.Range("A94").Value = 5750
$ws.Range("B94").Value = 45741.95833333334
$ws.Range("A95").Value = 5680
$ws.Range("B95").Value = 45741.96875
$ws.Range("A96").Value = 5640
$ws.Range("B96").Value = 45741.97916666666
$ws.Range("A97").Value = 5600
$ws.Range("B97").Value = 45741.98958333334
